# Add 2022-Q3 data:
#  - Insert a new "2022-Q3" sheet (copied from "2022-Q1" to inherit its layout/
#    styles) right after the "总计" summary sheet.
#  - Insert a new row into "总计" for the 2022-Q3 summary figures.

$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, $val) {
    # Force the cell to be stored as text even when the value looks numeric
    # (e.g. a fund code like "005632" or a decimal like "2.71"), matching the
    # inlineStr cells used throughout the workbook, then drop the now-unneeded
    # "@" number format so no stray style sticks to the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------------
# 1. Build the new "2022-Q3" worksheet by copying "2022-Q1" (same columns /
#    header / styling as every other quarterly sheet) and placing it right
#    after "总计".
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$wsQ1 = $wb.Worksheets.Item("2022-Q1")
$wsQ1.Copy($null, $wsTotal)

$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Name = "2022-Q3"

Set-TextValue $wsQ3.Cells.Item(2,2) "005632"
Set-TextValue $wsQ3.Cells.Item(2,3) "鹏华量化先锋混合"
Set-TextValue $wsQ3.Cells.Item(2,4) "2.71"
Set-TextValue $wsQ3.Cells.Item(2,5) "92.57"
Set-TextValue $wsQ3.Cells.Item(2,6) "1.81"
Set-TextValue $wsQ3.Cells.Item(2,7) "0.0491"
$wsQ3.Cells.Item(2,8).Value = 4

# ---------------------------------------------------------------------------
# 2. Insert the matching row into "总计", right above the existing 2022-Q1
#    row, pushing the older quarters down.
# ---------------------------------------------------------------------------
$wsTotal.Rows.Item(2).Insert()

$wsTotal.Cells.Item(2,2).Value = "2022-Q3"
$wsTotal.Cells.Item(2,3).Value = 1
$wsTotal.Cells.Item(2,4).Value = 0.05
$wsTotal.Range("B2:D2").ClearFormats()

# Re-number the leading index column (0,1,2,3,4) now that a row was inserted
# and copy the header-row-index cell style (s="2") down onto every row.
for ($r = 2; $r -le 6; $r++) {
    $wsTotal.Cells.Item($r,1).Value = $r - 2
}
$wsTotal.Cells.Item(3,1).Copy()
$wsTotal.Range("A2:A6").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Restore the originally-active tab (the last sheet, "2020-Q4").
# ---------------------------------------------------------------------------
$wsLast = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsLast.Activate()
